# Enhanced AssessmentCategoryType code table to test allegations dimension
# situation for Pima.
#
# Adds five new rows (AssessmentCategoryType 2 .. 6, ids 2..6) into the
# AssessmentCategoryType sheet, between the existing "1" row and the
# "99998"/"99999" sentinel rows, and leaves that sheet as the active
# sheet/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AssessmentCategoryType")

# Insert 5 blank rows right before the "99998" sentinel row (row 3),
# pushing the sentinel rows down to rows 8 and 9.
$ws.Rows("3:7").Insert()

# The inserted rows copied formatting from the row above (style "Normal 2");
# the new data rows should use the plain default style, matching the
# original "1" row once its style is cleared as well.
$ws.Range("A2:B7").Style = "Normal"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "AssessmentCategoryType 2"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "AssessmentCategoryType 3"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "AssessmentCategoryType 4"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "AssessmentCategoryType 5"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "AssessmentCategoryType 6"

# Make AssessmentCategoryType the active/selected sheet, with B12 selected,
# moving the "tabSelected" marker away from whichever sheet had it before.
$ws.Activate()
$ws.Range("B12").Select()
